## LoginInfo.xlsx edit
## - A1 label changes from "Login" to "Username" (B1/A2/B2 stay the same)
## - Active selection on Sheet1 moves from B2 to F11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in A1 (this also causes the shared-strings table
# to be rewritten: "Login" drops out and "Username" is appended at the end).
$ws.Range("A1").Value = "Username"

# Move / record the current selection on the sheet.
$ws.Range("F11").Select() | Out-Null
